$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.823 (0.815 ± 0.011)"
$ws.Range("C2").Value = "00:00:51 (00:01:03 ± 00:00:06)"
$ws.Range("D2").Value = "00:00:10 (00:00:10 ± 00:00:00)"
$ws.Range("B4").Value = "0.188 (0.116 ± 0.046)"
$ws.Range("C4").Value = "00:00:18 (00:00:25 ± 00:00:05)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B5").Value = "0.047 (0.011 ± 0.014)"
$ws.Range("C5").Value = "00:05:08 (00:05:23 ± 00:00:16)"
$ws.Range("D5").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B6").Value = "0.265 (0.182 ± 0.062)"
$ws.Range("C6").Value = "00:04:59 (00:05:03 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:01 (00:00:01 ± 00:00:00)"
$ws.Range("B9").Value = "0.272 (0.137 ± 0.061)"
$ws.Range("C9").Value = "00:05:00 (00:05:02 ± 00:00:02)"
$ws.Range("D9").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B11").Value = "0.290 (0.163 ± 0.069)"
$ws.Range("C11").Value = "00:05:11 (00:05:23 ± 00:00:10)"
$ws.Range("D11").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B13").Value = "0.026 (0.005 ± 0.009)"
$ws.Range("C13").Value = "00:00:01 (00:00:01 ± 00:00:00)"
$ws.Range("D13").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B14").Value = "0.220 (0.086 ± 0.059)"
$ws.Range("C14").Value = "00:04:15 (00:04:34 ± 00:00:14)"
$ws.Range("D14").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B15").Value = "0.201 (0.152 ± 0.050)"
$ws.Range("C15").Value = "00:00:55 (00:01:19 ± 00:00:42)"
$ws.Range("D15").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B16").Value = "0.273 (0.168 ± 0.057)"
$ws.Range("C16").Value = "00:00:49 (00:00:52 ± 00:00:01)"
$ws.Range("D16").Value = "00:00:00 (00:00:00 ± 00:00:00)"
